$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 99
$ws.Range("I33").Value = 99
$ws.Range("K33").Value = 99
$ws.Range("M33").Value = 130
$ws.Range("H132").Value = 807.8108
$ws.Range("I132").Value = 764.55884
$ws.Range("J132").Value = 1298
$ws.Range("K132").Value = 2293.67652
$ws.Range("L132").Value = 3894
$ws.Range("M132").Value = 236.32348
$ws.Range("N132").Value = -8954
$ws.Range("H138").Value = 1573
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1994.1595
$ws.Range("I32").Value = 1421.2561
$ws.Range("K32").Value = 1421.2561
$ws.Range("M32").Value = -1134.2561
$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30976
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").Value = $null
$ws.Range("H61").Value = 2924.4167
$ws.Range("I61").Value = 1959.75
$ws.Range("J61").Value = 4853.75
$ws.Range("K61").Value = 1959.75
$ws.Range("L61").Value = 4853.75
$ws.Range("M61").Value = -1747.75
$ws.Range("N61").Value = -5277.75
$ws.Range("H74").Value = 1487.4445
$ws.Range("I74").Value = 603.9167
$ws.Range("J74").Value = 3254.5
$ws.Range("K74").Value = 603.9167
$ws.Range("L74").Value = 3254.5
$ws.Range("M74").Value = 270.0833
$ws.Range("N74").Value = -5002.5
$ws.Range("H77").Value = 1487.4445
$ws.Range("I77").Value = 603.9167
$ws.Range("J77").Value = 3254.5
$ws.Range("K77").Value = 3019.5835
$ws.Range("L77").Value = 16272.5
$ws.Range("M77").Value = 1348.4165
$ws.Range("N77").Value = -25008.5
$ws.Range("H110").Value = 740.0833
$ws.Range("I110").Value = 740.0833
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 740.0833
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1304.9167
$ws.Range("N110").Value = $null
$ws.Range("H122").Value = 96719.75
$ws.Range("I122").Value = 128126.336
$ws.Range("K122").Value = 384379.008
$ws.Range("M122").Value = -381929.008
$ws.Range("H132").Value = 2675.577
$ws.Range("I132").Value = 2481.9583
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 7445.874899999999
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -4915.874899999999
$ws.Range("N132").Value = -20057
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null
$ws.Range("H135").Value = 100429
$ws.Range("J135").Value = 100429
$ws.Range("L135").Value = 100429
$ws.Range("N135").Value = -110569
$ws.Range("H136").Value = 2924.4167
$ws.Range("I136").Value = 1959.75
$ws.Range("J136").Value = 4853.75
$ws.Range("K136").Value = 5879.25
$ws.Range("L136").Value = 14561.25
$ws.Range("M136").Value = -3329.25
$ws.Range("N136").Value = -19661.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 500726.25
$ws.Range("I86").Value = 952.5
$ws.Range("J86").Value = 1000500
$ws.Range("K86").Value = 952.5
$ws.Range("L86").Value = 1000500
$ws.Range("M86").Value = 170.5
$ws.Range("N86").Value = -1002746
$ws.Range("H89").Value = 500726.25
$ws.Range("I89").Value = 952.5
$ws.Range("J89").Value = 1000500
$ws.Range("K89").Value = 4762.5
$ws.Range("L89").Value = 5002500
$ws.Range("M89").Value = 853.5
$ws.Range("N89").Value = -5013732
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1480.4667
$ws.Range("I31").Value = 948.8333
$ws.Range("J31").Value = 1834.8889
$ws.Range("K31").Value = 948.8333
$ws.Range("L31").Value = 1834.8889
$ws.Range("M31").Value = -653.8333
$ws.Range("N31").Value = -2424.8889
$ws.Range("H34").Value = 1480.4667
$ws.Range("I34").Value = 948.8333
$ws.Range("J34").Value = 1834.8889
$ws.Range("K34").Value = 948.8333
$ws.Range("L34").Value = 1834.8889
$ws.Range("M34").Value = -746.8333
$ws.Range("N34").Value = -2238.8889
$ws.Range("H58").Value = 2289679.5
$ws.Range("I58").Value = 3953977.8
$ws.Range("J58").Value = 1269.25
$ws.Range("K58").Value = 3953977.8
$ws.Range("L58").Value = 1269.25
$ws.Range("M58").Value = -3953774.8
$ws.Range("N58").Value = -1675.25
$ws.Range("H62").Value = 1666.6666
$ws.Range("I62").Value = 1500
$ws.Range("K62").Value = 1500
$ws.Range("M62").Value = -876
$ws.Range("H65").Value = 1666.6666
$ws.Range("I65").Value = 1500
$ws.Range("K65").Value = 7500
$ws.Range("M65").Value = -4380
$ws.Range("H94").Value = 1500
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1500
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1500
$ws.Range("M94").Value = $null
$ws.Range("N94").Value = -2402
$ws.Range("H136").Value = 2289679.5
$ws.Range("I136").Value = 3953977.8
$ws.Range("J136").Value = 1269.25
$ws.Range("K136").Value = 11861933.4
$ws.Range("L136").Value = 3807.75
$ws.Range("M136").Value = -11859383.4
$ws.Range("N136").Value = -8907.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 2619.4
$ws.Range("I7").Value = 3333.3333
$ws.Range("K7").Value = 9999.999899999999
$ws.Range("M7").Value = -9887.999899999999
$ws.Range("H80").Value = 3799.6924
$ws.Range("I80").Value = 4982.8335
$ws.Range("J80").Value = 2785.5715
$ws.Range("K80").Value = 14948.5005
$ws.Range("L80").Value = 8356.7145
$ws.Range("M80").Value = -14012.5005
$ws.Range("N80").Value = -10228.7145
$ws.Range("H83").Value = 3799.6924
$ws.Range("I83").Value = 4982.8335
$ws.Range("J83").Value = 2785.5715
$ws.Range("K83").Value = 44845.5015
$ws.Range("L83").Value = 25070.1435
$ws.Range("M83").Value = -40165.5015
$ws.Range("N83").Value = -34430.1435
$ws.Range("H92").Value = 807.5
$ws.Range("J92").Value = 880
$ws.Range("L92").Value = 2640
$ws.Range("N92").Value = -5136
$ws.Range("H129").Value = 45759.938
$ws.Range("I129").Value = 559.1
$ws.Range("J129").Value = 121094.664
$ws.Range("K129").Value = 1677.3
$ws.Range("L129").Value = 363283.992
$ws.Range("M129").Value = 3322.7
$ws.Range("N129").Value = -373283.992
$ws.Range("H131").Value = 11130139
$ws.Range("I131").Value = 100000350
$ws.Range("J131").Value = 21362.3
$ws.Range("K131").Value = 300001050
$ws.Range("L131").Value = 64086.89999999999
$ws.Range("M131").Value = -299996010
$ws.Range("N131").Value = -74166.89999999999
$ws.Range("H137").Value = 4034.8462
$ws.Range("J137").Value = 5914.7144
$ws.Range("L137").Value = 17744.1432
$ws.Range("N137").Value = -27944.1432
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 297
$ws.Range("I2").Value = 342.75
$ws.Range("J2").Value = 205.5
$ws.Range("K2").Value = 342.75
$ws.Range("L2").Value = 205.5
$ws.Range("M2").Value = -229.75
$ws.Range("N2").Value = -431.5
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = $null
$ws.Range("H132").Value = 1426881.6
$ws.Range("I132").Value = 1924954.8
$ws.Range("J132").Value = 3815.5715
$ws.Range("K132").Value = 5774864.4
$ws.Range("L132").Value = 11446.7145
$ws.Range("M132").Value = -5772334.4
$ws.Range("N132").Value = -16506.7145
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4025.7144
$ws.Range("I22").Value = 950
$ws.Range("J22").Value = 5256
$ws.Range("K22").Value = 950
$ws.Range("L22").Value = 5256
$ws.Range("M22").Value = -655
$ws.Range("N22").Value = -5846
$ws.Range("H27").Value = 4025.7144
$ws.Range("I27").Value = 950
$ws.Range("J27").Value = 5256
$ws.Range("K27").Value = 950
$ws.Range("L27").Value = 5256
$ws.Range("M27").Value = -843
$ws.Range("N27").Value = -5470
$ws.Range("H68").Value = 3828.3333
$ws.Range("I68").Value = 3594
$ws.Range("K68").Value = 3594
$ws.Range("M68").Value = -2845
$ws.Range("H71").Value = 3828.3333
$ws.Range("I71").Value = 3594
$ws.Range("K71").Value = 17970
$ws.Range("M71").Value = -14226
$ws.Range("H82").Value = 2570.5
$ws.Range("I82").Value = 996
$ws.Range("J82").Value = 3357.75
$ws.Range("K82").Value = 996
$ws.Range("L82").Value = 3357.75
$ws.Range("M82").Value = -635
$ws.Range("N82").Value = -4079.75
$ws.Range("H85").Value = 2570.5
$ws.Range("I85").Value = 996
$ws.Range("J85").Value = 3357.75
$ws.Range("K85").Value = 996
$ws.Range("L85").Value = 3357.75
$ws.Range("M85").Value = 252
$ws.Range("N85").Value = -5853.75
$ws.Range("H98").Value = 49500
$ws.Range("J98").Value = 49500
$ws.Range("L98").Value = 49500
$ws.Range("N98").Value = -55490
$ws.Range("H122").Value = 11874.625
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550
$ws.Range("H136").Value = 3680.3928
$ws.Range("I136").Value = 2345.5789
$ws.Range("J136").Value = 6498.3335
$ws.Range("K136").Value = 7036.736699999999
$ws.Range("L136").Value = 19495.0005
$ws.Range("M136").Value = -4486.736699999999
$ws.Range("N136").Value = -24595.0005
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4175
$ws.Range("J96").Value = 4175
$ws.Range("L96").Value = 4175
$ws.Range("N96").Value = -6921
$ws.Range("H101").Value = 19999
$ws.Range("J101").Value = 19999
$ws.Range("L101").Value = 19999
$ws.Range("N101").Value = -26489
$ws.Range("H128").Value = 30000
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960
$ws.Range("H136").Value = 11823572
$ws.Range("I136").Value = 18521520
$ws.Range("J136").Value = 3663.8235
$ws.Range("K136").Value = 55564560
$ws.Range("L136").Value = 10991.4705
$ws.Range("M136").Value = -55562010
$ws.Range("N136").Value = -16091.4705
